$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the cells we touch so numeric-looking strings
# (prices, percentages) are preserved as Text, matching the source data.
$textCells = @("D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "B7", "C7", "D7", "E7", "B8", "C8", "D8", "E8", "B9", "C9", "D9", "E9", "B10", "C10", "D10", "E10", "B11", "C11", "D11", "E11", "B12", "C12", "D12", "E12", "B13", "C13", "D13", "E13", "B14", "C14", "D14", "E14", "B15", "C15", "D15", "E15", "B16", "C16", "D16", "E16", "B17", "C17", "D17", "E17", "E18", "E19", "D20", "E20", "D21", "E21", "D23", "E23", "D24", "E24", "E25", "D27", "E27", "D39", "E39", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "D46", "E46", "D47", "D48", "D49", "E49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values
$ws.Range("D2").Value = "308.74"
$ws.Range("E2").Value = "1.15%"
$ws.Range("D3").Value = "38.64"
$ws.Range("E3").Value = "8.05%"
$ws.Range("D4").Value = "5.101"
$ws.Range("E4").Value = "1.25%"
$ws.Range("D5").Value = "0.08117"
$ws.Range("E5").Value = "1.04%"
$ws.Range("D6").Value = "1.969"
$ws.Range("E6").Value = "5.18%"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").Value = "4.209"
$ws.Range("E7").Value = "1.56%"
$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D8").Value = "7.940"
$ws.Range("E8").Value = "1.90%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "0.9279"
$ws.Range("E9").Value = "0.90%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "0.1422"
$ws.Range("E10").Value = "11.97%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "0.1962"
$ws.Range("E11").Value = "2.45%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "0.09134"
$ws.Range("E12").Value = "0.53%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "0.03506"
$ws.Range("E13").Value = "1.56%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "0.09809"
$ws.Range("E14").Value = "-0.48%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "0.001408"
$ws.Range("E15").Value = "0.21%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "0.006079"
$ws.Range("E16").Value = "-1.36%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "3.661"
$ws.Range("E17").Value = "-4.13%"
$ws.Range("E18").Value = "2.14%"
$ws.Range("E19").Value = "1.29%"
$ws.Range("D20").Value = "0.1303"
$ws.Range("E20").Value = "-1.33%"
$ws.Range("D21").Value = "4.805"
$ws.Range("E21").Value = "-8.25%"
$ws.Range("D23").Value = "0.04430"
$ws.Range("E23").Value = "0.05%"
$ws.Range("D24").Value = "0.001218"
$ws.Range("E24").Value = "-1.27%"
$ws.Range("E25").Value = "4.74%"
$ws.Range("D27").Value = "0.0001301"
$ws.Range("E27").Value = "3.96%"
$ws.Range("D39").Value = "0.02098"
$ws.Range("E39").Value = "7.86%"
$ws.Range("D40").Value = "0.05145"
$ws.Range("E40").Value = "-2.66%"
$ws.Range("D41").Value = "0.007480"
$ws.Range("E41").Value = "-2.02%"
$ws.Range("D42").Value = "0.01011"
$ws.Range("E42").Value = "-0.30%"
$ws.Range("D43").Value = "0.1359"
$ws.Range("E43").Value = "0.32%"
$ws.Range("D44").Value = "0.002142"
$ws.Range("E44").Value = "-1.42%"
$ws.Range("D45").Value = "0.009226"
$ws.Range("E45").Value = "-4.14%"
$ws.Range("D46").Value = "0.00006374"
$ws.Range("E46").Value = "4.44%"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("D48").Value = "0.003068"
$ws.Range("D49").Value = "0.001601"
$ws.Range("E49").Value = "-3.57%"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("D51").Value = "0.0002002"
